{"js": "const replacements = [\n  [\"2025-02-08 Saturday\", \"2025-02-09 Sunday\"],\n  [\"185\u00f75=\", \"140\u00f75=\"],\n  [\"449\u00f79=\", \"628\u00f74=\"],\n  [\"840\u00f79=\", \"987\u00f73=\"],\n  [\"505\u00f78=\", \"812\u00f79=\"],\n  [\"289\u00f72=\", \"478\u00f73=\"],\n  [\"762\u00f76=\", \"123\u00f78=\"],\n  [\"643\u00f78=\", \"264\u00f78=\"],\n  [\"322\u00f77=\", \"524\u00f76=\"],\n  [\"597\u00f72=\", \"218\u00f75=\"],\n  [\"251\u00f77=\", \"434\u00f79=\"],\n  [\"427\u00f73=\", \"955\u00f75=\"],\n  [\"708\u00f73=\", \"653\u00f79=\"],\n  [\"373\u00f76=\", \"918\u00f75=\"],\n  [\"442\u00f73=\", \"428\u00f78=\"],\n  [\"330\u00f72=\", \"290\u00f72=\"],\n  [\"844\u00f74=\", \"686\u00f77=\"],\n  [\"202\u00f79=\", \"904\u00f74=\"],\n  [\"641\u00f75=\", \"883\u00f72=\"],\n  [\"824\u00f75=\", \"509\u00f78=\"],\n  [\"226\u00f76=\", \"607\u00f74=\"],\n  [\"173\u00f72=\", \"890\u00f74=\"],\n  [\"774\u00f76=\", \"873\u00f72=\"],\n  [\"785\u00f74=\", \"840\u00f75=\"],\n  [\"845\u00f73=\", \"291\u00f78=\"],\n  [\"683\u00f74=\", \"565\u00f76=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-08 Saturday\", \"2025-02-09 Sunday\"),\n    @(\"185\u00f75=\", \"140\u00f75=\"),\n    @(\"449\u00f79=\", \"628\u00f74=\"),\n    @(\"840\u00f79=\", \"987\u00f73=\"),\n    @(\"505\u00f78=\", \"812\u00f79=\"),\n    @(\"289\u00f72=\", \"478\u00f73=\"),\n    @(\"762\u00f76=\", \"123\u00f78=\"),\n    @(\"643\u00f78=\", \"264\u00f78=\"),\n    @(\"322\u00f77=\", \"524\u00f76=\"),\n    @(\"597\u00f72=\", \"218\u00f75=\"),\n    @(\"251\u00f77=\", \"434\u00f79=\"),\n    @(\"427\u00f73=\", \"955\u00f75=\"),\n    @(\"708\u00f73=\", \"653\u00f79=\"),\n    @(\"373\u00f76=\", \"918\u00f75=\"),\n    @(\"442\u00f73=\", \"428\u00f78=\"),\n    @(\"330\u00f72=\", \"290\u00f72=\"),\n    @(\"844\u00f74=\", \"686\u00f77=\"),\n    @(\"202\u00f79=\", \"904\u00f74=\"),\n    @(\"641\u00f75=\", \"883\u00f72=\"),\n    @(\"824\u00f75=\", \"509\u00f78=\"),\n    @(\"226\u00f76=\", \"607\u00f74=\"),\n    @(\"173\u00f72=\", \"890\u00f74=\"),\n    @(\"774\u00f76=\", \"873\u00f72=\"),\n    @(\"785\u00f74=\", \"840\u00f75=\"),\n    @(\"845\u00f73=\", \"291\u00f78=\"),\n    @(\"683\u00f74=\", \"565\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n}\n"}
